$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column X (a new attendance-taking session / date column)
#    before the existing "Faltas" (X, now shifted to Y) summary column.
#    Excel's native column-insert already takes care of:
#      - shifting the old X column (Faltas) data + formulas to Y
#      - de-sharing the shared formula into individual COUNTIF(...) formulas
#        that keep referencing I:W (NOT extended to include the new column)
#      - updating dimension / row spans
#      - the new column inherits formatting (style) from the column to its
#        left (W), which is what we want for the new attendance column
# ---------------------------------------------------------------------------
$ws.Columns("X").Insert()

# ---------------------------------------------------------------------------
# 2. New session date in X4 (26-Jun-2023 == serial 45103)
# ---------------------------------------------------------------------------
$ws.Range("X4").Value = 45103

# ---------------------------------------------------------------------------
# 3. Fill in attendance marks ("F") for the new column - copied from column W
#    for most students, except rows 29 and 34 where the monitor left the new
#    column blank (row 29's W mark was also corrected from "F" to blank).
# ---------------------------------------------------------------------------
$rowsWithF = @(5,10,11,13,15,16,18,19,20,21,22,24,31,32,33,35,36,37)
foreach ($r in $rowsWithF) {
    $ws.Cells.Item($r, 24).Value = "F"
}

# Row 29 correction: clear the old "F" mark in column W (col 23)
$ws.Cells.Item(29, 23).Value = ""

# ---------------------------------------------------------------------------
# 4. Conditional formatting bookkeeping.
#    After the column insert, the CF rule that used to target the Faltas
#    column (X5:X37, ">4") now needs to target Y5:Y37 instead (it kept
#    pointing at X, which is now the blank new column).
# ---------------------------------------------------------------------------
$faltasRule = $ws.Range("X5:X37").FormatConditions.Item(1)
$faltasRule.ModifyAppliesToRange($ws.Range("Y5:Y37"))

# The rules that used to apply only to V5:V37 get extended to also cover the
# new column (so they now read V5:W37), mirroring the column that was
# inserted next to them during the edit.
$vRange = $ws.Range("V5:V37")
$vCount = $vRange.FormatConditions.Count
for ($i = 1; $i -le $vCount; $i++) {
    $rule = $vRange.FormatConditions.Item($i)
    $rule.ModifyAppliesToRange($ws.Range("V5:W37"))
}

# ---------------------------------------------------------------------------
# 5. Re-create, on the brand-new X5:X37 column, the full set of attendance
#    conditional-formatting rules (the same rules that already exist on
#    W5:W37 and used to exist on V5:V37), each with its own dxf (font +
#    fill) so the new column is highlighted exactly like its neighbours.
# ---------------------------------------------------------------------------
$xRange = $ws.Range("X5:X37")

function Add-CfRule($range, $operator, $formula, $fontColor, $fillColor) {
    $rule = $range.FormatConditions.Add(1, $operator, $formula)
    $rule.Font.Color = $fontColor
    $rule.Interior.Color = $fillColor
    return $rule
}

$RED_FONT = 393372
$RED_FILL = 13551615
$GREEN_FONT = 24832
$GREEN_FILL = 13561798

# Copy of the W5:W37 rule set
Add-CfRule $xRange 3 '"F"' $RED_FONT $RED_FILL | Out-Null
Add-CfRule $xRange 6 '1' $GREEN_FONT $GREEN_FILL | Out-Null
Add-CfRule $xRange 6 '0' $RED_FONT $RED_FILL | Out-Null
Add-CfRule $xRange 3 '"F"' $RED_FONT $RED_FILL | Out-Null
Add-CfRule $xRange 3 '0' $GREEN_FONT $GREEN_FILL | Out-Null
Add-CfRule $xRange 3 '" "' $GREEN_FONT $GREEN_FILL | Out-Null
Add-CfRule $xRange 3 '"F"' $RED_FONT $RED_FILL | Out-Null

# Copy of the V5:V37 rule set
Add-CfRule $xRange 3 '"F"' $RED_FONT $RED_FILL | Out-Null
Add-CfRule $xRange 6 '1' $GREEN_FONT $GREEN_FILL | Out-Null
Add-CfRule $xRange 6 '0' $RED_FONT $RED_FILL | Out-Null
Add-CfRule $xRange 3 '"F"' $RED_FONT $RED_FILL | Out-Null
Add-CfRule $xRange 3 '0' $GREEN_FONT $GREEN_FILL | Out-Null
Add-CfRule $xRange 3 '" "' $GREEN_FONT $GREEN_FILL | Out-Null
Add-CfRule $xRange 3 '"F"' $RED_FONT $RED_FILL | Out-Null

# ---------------------------------------------------------------------------
# 6. Restore selection on the sheet (user ended up with X34 selected in the
#    frozen right-hand pane).
# ---------------------------------------------------------------------------
$ws.Range("X34").Select()
